# Auto-generated cell value updates for Asura_Profits workbook
# (scheduled runner refresh of currentAveragePrice / LevePrice / LeveProfit columns)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 254.77777
$ws.Range("I33").Value = 277.56522
$ws.Range("K33").Value = 277.56522
$ws.Range("M33").Value = -48.56522000000001
$ws.Range("H51").Value = 3241.4285
$ws.Range("I51").Value = 2747.5
$ws.Range("J51").Value = 3439
$ws.Range("K51").Value = 2747.5
$ws.Range("L51").Value = 3439
$ws.Range("M51").Value = -2263.5
$ws.Range("N51").Value = -4407
$ws.Range("H129").Value = 1128.0984
$ws.Range("J129").Value = 1167.4912
$ws.Range("L129").Value = 3502.4736
$ws.Range("N129").Value = -13502.4736
$ws.Range("H132").Value = 1950.66
$ws.Range("I132").Value = 1846.6857
$ws.Range("J132").Value = 2193.2666
$ws.Range("K132").Value = 5540.0571
$ws.Range("L132").Value = 6579.7998
$ws.Range("M132").Value = -3010.0571
$ws.Range("N132").Value = -11639.7998

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 8585.5
$ws.Range("I31").Value = 8585.5
$ws.Range("K31").Value = 8585.5
$ws.Range("M31").Value = -8291.5
$ws.Range("H32").Value = 14029.981
$ws.Range("I32").Value = 16604.953
$ws.Range("J32").Value = 3964.182
$ws.Range("K32").Value = 16604.953
$ws.Range("L32").Value = 3964.182
$ws.Range("M32").Value = -16317.953
$ws.Range("N32").Value = -4538.182
$ws.Range("H45").Value = 1367.7273
$ws.Range("I45").Value = 1345.8572
$ws.Range("J45").Value = 1406
$ws.Range("K45").Value = 1345.8572
$ws.Range("L45").Value = 1406
$ws.Range("M45").Value = -968.8571999999999
$ws.Range("N45").Value = -2160
$ws.Range("H109").Value = 34833
$ws.Range("J109").Value = 34833
$ws.Range("L109").Value = 34833
$ws.Range("N109").Value = -37607
$ws.Range("H122").Value = 3737.6538
$ws.Range("I122").Value = 3318.9285
$ws.Range("J122").Value = 4226.1665
$ws.Range("K122").Value = 9956.7855
$ws.Range("L122").Value = 12678.4995
$ws.Range("M122").Value = -7506.7855
$ws.Range("N122").Value = -17578.4995
$ws.Range("H132").Value = 1912.6471
$ws.Range("I132").Value = 1257.3226
$ws.Range("J132").Value = 2928.4
$ws.Range("K132").Value = 3771.9678
$ws.Range("L132").Value = 8785.200000000001
$ws.Range("M132").Value = -1241.9678
$ws.Range("N132").Value = -13845.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 27854.428
$ws.Range("J108").Value = 27854.428
$ws.Range("L108").Value = 27854.428
$ws.Range("N108").Value = -35534.428
$ws.Range("H134").Value = 2350
$ws.Range("I134").Value = 1905
$ws.Range("J134").Value = 2604.2856
$ws.Range("K134").Value = 5715
$ws.Range("L134").Value = 7812.8568
$ws.Range("M134").Value = -3180
$ws.Range("N134").Value = -12882.8568

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1685.6086
$ws.Range("I31").Value = 1426.4706
$ws.Range("J31").Value = 2419.8333
$ws.Range("K31").Value = 1426.4706
$ws.Range("L31").Value = 2419.8333
$ws.Range("M31").Value = -1131.4706
$ws.Range("N31").Value = -3009.8333
$ws.Range("H34").Value = 1685.6086
$ws.Range("I34").Value = 1426.4706
$ws.Range("J34").Value = 2419.8333
$ws.Range("K34").Value = 1426.4706
$ws.Range("L34").Value = 2419.8333
$ws.Range("M34").Value = -1224.4706
$ws.Range("N34").Value = -2823.8333

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 10000
$ws.Range("I110").Value = 5000
$ws.Range("J110").Value = 15000
$ws.Range("K110").Value = 15000
$ws.Range("L110").Value = 45000
$ws.Range("M110").Value = -10910
$ws.Range("N110").Value = -53180
$ws.Range("H113").Value = 641.14703
$ws.Range("I113").Value = 600
$ws.Range("J113").Value = 693.26666
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 2079.79998
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -6419.79998

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5619.95
$ws.Range("I70").Value = 5506.6
$ws.Range("J70").Value = 5960
$ws.Range("K70").Value = 5506.6
$ws.Range("L70").Value = 5960
$ws.Range("M70").Value = -5236.6
$ws.Range("N70").Value = -6500
$ws.Range("H73").Value = 5619.95
$ws.Range("I73").Value = 5506.6
$ws.Range("J73").Value = 5960
$ws.Range("K73").Value = 5506.6
$ws.Range("L73").Value = 5960
$ws.Range("M73").Value = -4570.6
$ws.Range("N73").Value = -7832
$ws.Range("H102").Value = 3408.5386
$ws.Range("I102").Value = 3580
$ws.Range("J102").Value = 3301.375
$ws.Range("K102").Value = 3580
$ws.Range("L102").Value = 3301.375
$ws.Range("M102").Value = -1958
$ws.Range("N102").Value = -6545.375
$ws.Range("H126").Value = 3773
$ws.Range("I126").Value = 3669.5
$ws.Range("J126").Value = 3980
$ws.Range("K126").Value = 11008.5
$ws.Range("L126").Value = 11940
$ws.Range("M126").Value = -8538.5
$ws.Range("N126").Value = -16880

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6288.3335
$ws.Range("I40").Value = 6746
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 6746
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -6610
$ws.Range("N40").Value = -4272
$ws.Range("H122").Value = 12504406
$ws.Range("I122").Value = 4374.8335
$ws.Range("K122").Value = 13124.5005
$ws.Range("M122").Value = -10674.5005
$ws.Range("H132").Value = 5777.625
$ws.Range("I132").Value = 6610.6665
$ws.Range("J132").Value = 4389.222
$ws.Range("K132").Value = 19831.9995
$ws.Range("L132").Value = 13167.666
$ws.Range("M132").Value = -17301.9995
$ws.Range("N132").Value = -18227.666

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 76666.664
$ws.Range("H73").Value = 76666.664
$ws.Range("H118").Value = 31399.334
$ws.Range("J118").Value = 31399.334
$ws.Range("L118").Value = 31399.334
$ws.Range("N118").Value = -34713.334
$ws.Range("H121").Value = 26182.117
$ws.Range("J121").Value = 26182.117
$ws.Range("L121").Value = 26182.117
$ws.Range("N121").Value = -29676.117
$ws.Range("H123").Value = 40429
$ws.Range("J123").Value = 40429
$ws.Range("L123").Value = 40429
$ws.Range("N123").Value = -50229
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H127").Value = 40429
$ws.Range("J127").Value = 40429
$ws.Range("L127").Value = 40429
$ws.Range("N127").Value = -50349
$ws.Range("H132").Value = 2724.7778
$ws.Range("I132").Value = 2143
$ws.Range("J132").Value = 3015.6667
$ws.Range("K132").Value = 6429
$ws.Range("L132").Value = 9047.000100000001
$ws.Range("M132").Value = -3899
$ws.Range("N132").Value = -14107.0001
$ws.Range("H136").Value = 1611.1305
$ws.Range("I136").Value = 1452.85
$ws.Range("K136").Value = 4358.549999999999
$ws.Range("M136").Value = -1808.549999999999
